$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.702.06"
$ws.Range("E2").Value = "  +0.27%  "

# Row 3
$ws.Range("D3").Value = "1.600.12"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("E5").Value = "  +0.13%  "

# Row 6
$ws.Range("E6").Value = "  -0.80%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.43%  "

# Row 9
$ws.Range("E9").Value = "  +1.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.87%  "

# Row 12
$ws.Range("D12").Value = "1.824.52"
$ws.Range("E12").Value = "  +0.15%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.584.89"
$ws.Range("E13").Value = "  +1.18%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.65%  "

# Row 15
$ws.Range("E15").Value = "  +0.55%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.01%  "

# Row 17
$ws.Range("D17").Value = "26.677.08"

# Row 18
$ws.Range("E18").Value = "  +4.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.22"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.37%  "

# Row 21
$ws.Range("E21").Value = "  +0.16%  "

# Row 22
$ws.Range("E22").Value = "  +0.69%  "

# Row 23
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.94"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.12"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.76%  "

# Row 26
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("E27").Value = "  -0.03%  "

# Row 28
$ws.Range("E28").Value = "  +0.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.96%  "

# Row 31
$ws.Range("E31").Value = "  +0.01%  "

# Row 32
$ws.Range("E32").Value = "  +0.77%  "

# Row 33
$ws.Range("E33").Value = "  +1.93%  "

# Row 34
$ws.Range("D34").Value = "1.290.89"
$ws.Range("E34").Value = "  +0.77%  "

# Row 35
$ws.Range("E35").Value = "  -4.89%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.47"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.89%  "

# Row 37
$ws.Range("E37").Value = "  +0.42%  "

# Row 38
$ws.Range("E38").Value = "  -0.07%  "

# Row 39
$ws.Range("E39").Value = "  +16.07%  "

# Row 40
$ws.Range("E40").Value = "  -1.87%  "

# Row 41
$ws.Range("E41").Value = "  -0.31%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("E43").Value = "  -0.72%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.06%  "

# Row 45
$ws.Range("D45").Value = "1.735.46"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.92%  "

# Row 47
$ws.Range("E47").Value = "  -1.19%  "

# Row 48
$ws.Range("E48").Value = "  -1.00%  "

# Row 49
$ws.Range("E49").Value = "  +0.71%  "

# Row 50
$ws.Range("E50").Value = "  +0.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.39"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.63%  "
